# CampaignEvaluationTool.xlsx - add basic ROI calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells S1:U1 ("Costs", "Extra Revenue", "ROI")
# ---------------------------------------------------------------------------
$ws.Range("S1").Value = "Costs"
$ws.Range("T1").Value = "Extra Revenue"
$ws.Range("U1").Value = "ROI"
$ws.Range("S1:U1").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Row 4 - rework the existing S4 formula, add T4/U4
# ---------------------------------------------------------------------------
$ws.Range("S4").Formula = "=P4*0.8"
$ws.Range("T4").Formula = "=10*24*P4*C4"
$ws.Range("U4").Formula = "=T4/S4"
$ws.Range("S4:U4").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Rows 5 and 6 - new campaign scenarios, copied from row 4's layout/style
# ---------------------------------------------------------------------------
$ws.Range("A4:R4").Copy($ws.Range("A5:R5"))
$ws.Range("A4:R4").Copy($ws.Range("A6:R6"))

# --- Row 5 inputs ---
$ws.Range("A5").Value = 0.009
$ws.Range("B5").Value = 0.005
$ws.Range("E5").Value = 0.05
$ws.Range("F5").Value = 0.8
$ws.Range("N5").Value = 0.75
$ws.Range("O5").Value = 30000

# --- Row 5 formulas (restore formula cells overwritten by the value copy) ---
$ws.Range("C5").Formula = "=A5-B5"
$ws.Range("D5").Formula = "=A5/B5-1"
$ws.Range("G5").Formula = "=(A5*P5+B5*Q5)/O5"
$ws.Range("H5").Formula = "=NORM.S.INV(1-E5)"
$ws.Range("I5").Formula = "=NORM.S.INV(F5)"
$ws.Range("J5").Formula = "=SQRT(G5*(1-G5)*(1/Q5+1/P5))"
$ws.Range("K5").Formula = "=(A5-B5)/J5"
$ws.Range("L5").Formula = "=H5+I5"
$ws.Range("M5").Formula = "=(NORM.S.DIST(K5,TRUE))"
$ws.Range("P5").Formula = "=O5*N5"
$ws.Range("Q5").Formula = "=O5-P5"
$ws.Range("R5").Formula = "=IF(K5>L5,""Yes"",""No"")"

# --- Row 5 new ROI columns ---
$ws.Range("S5").Formula = "=P5*0.8"
$ws.Range("T5").Formula = "=10*24*P5*C5"
$ws.Range("U5").Formula = "=T5/S5"
$ws.Range("S5:U5").Style = "Normal"

# --- Row 6 inputs ---
$ws.Range("A6").Value = 0.012
$ws.Range("B6").Value = 0.008
$ws.Range("E6").Value = 0.05
$ws.Range("F6").Value = 0.8
$ws.Range("N6").Value = 0.55
$ws.Range("O6").Value = 30000

# --- Row 6 formulas (restore formula cells overwritten by the value copy) ---
$ws.Range("C6").Formula = "=A6-B6"
$ws.Range("D6").Formula = "=A6/B6-1"
$ws.Range("G6").Formula = "=(A6*P6+B6*Q6)/O6"
$ws.Range("H6").Formula = "=NORM.S.INV(1-E6)"
$ws.Range("I6").Formula = "=NORM.S.INV(F6)"
$ws.Range("J6").Formula = "=SQRT(G6*(1-G6)*(1/Q6+1/P6))"
$ws.Range("K6").Formula = "=(A6-B6)/J6"
$ws.Range("L6").Formula = "=H6+I6"
$ws.Range("M6").Formula = "=(NORM.S.DIST(K6,TRUE))"
$ws.Range("P6").Formula = "=O6*N6"
$ws.Range("Q6").Formula = "=O6-P6"
$ws.Range("R6").Formula = "=IF(K6>L6,""Yes"",""No"")"

# --- Row 6 new ROI columns ---
$ws.Range("S6").Formula = "=P6*0.8"
$ws.Range("T6").Formula = "=10*24*P6*C6"
$ws.Range("U6").Formula = "=T6/S6"
$ws.Range("S6:U6").Style = "Normal"

# ---------------------------------------------------------------------------
# 4. View tidy-up: scroll back to A1 and select U4 (matches the saved file)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("U4").Select()
